# Auto-generated Excel COM-interop script
# Fixes: (1) Absentees bug in "Student Summary" sheet (DSPC608 / DSPC607 stats + rounding)
#        (2) Consolidated output bug in "Slow Learners" / "Fast Learners" sheets (re-ranked rosters)

$wb = $excel.ActiveWorkbook

# ---- Sheet 1: Student Summary ----
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A1").Value = "Course Code"
$ws1.Range("B1").Value = "DSPC604"
$ws1.Range("C1").Value = "DSPE605"
$ws1.Range("D1").Value = "DSPC608"
$ws1.Range("E1").Value = "DSPC601"
$ws1.Range("F1").Value = "EEOE 606"
$ws1.Range("G1").Value = "DSPC602"
$ws1.Range("H1").Value = "DSPE603"
$ws1.Range("I1").Value = "DSPC607"

$ws1.Range("A2").Value = "Total Students"
$ws1.Range("B2").Value = 51
$ws1.Range("C2").Value = 51
$ws1.Range("D2").Value = 51
$ws1.Range("E2").Value = 51
$ws1.Range("F2").Value = 51
$ws1.Range("G2").Value = 51
$ws1.Range("H2").Value = 51
$ws1.Range("I2").Value = 51

$ws1.Range("A3").Value = "Total Students Appeared"
$ws1.Range("B3").Value = 51
$ws1.Range("C3").Value = 51
$ws1.Range("D3").Value = 49
$ws1.Range("E3").Value = 51
$ws1.Range("F3").Value = 51
$ws1.Range("G3").Value = 51
$ws1.Range("H3").Value = 51
$ws1.Range("I3").Value = 48

$ws1.Range("A4").Value = "Total Absent"
$ws1.Range("B4").Value = 0
$ws1.Range("C4").Value = 0
$ws1.Range("D4").Value = 2
$ws1.Range("E4").Value = 0
$ws1.Range("F4").Value = 0
$ws1.Range("G4").Value = 0
$ws1.Range("H4").Value = 0
$ws1.Range("I4").Value = 3

$ws1.Range("A5").Value = "Average Marks"
$ws1.Range("B5").Value = 12.88
$ws1.Range("C5").Value = 19.04
$ws1.Range("D5").Value = 7.12
$ws1.Range("E5").Value = 26.84
$ws1.Range("F5").Value = 22.22
$ws1.Range("G5").Value = 19.88
$ws1.Range("H5").Value = 27.24
$ws1.Range("I5").Value = 11.71

$ws1.Range("A6").Value = "Students Less than 40%"
$ws1.Range("B6").Value = 35
$ws1.Range("C6").Value = 15
$ws1.Range("D6").Value = 28
$ws1.Range("E6").Value = 6
$ws1.Range("F6").Value = 11
$ws1.Range("G6").Value = 15
$ws1.Range("H6").Value = 7
$ws1.Range("I6").Value = 3

$ws1.Range("A7").Value = "Students Between 40 %. and 75 %"
$ws1.Range("B7").Value = 16
$ws1.Range("C7").Value = 31
$ws1.Range("D7").Value = 4
$ws1.Range("E7").Value = 22
$ws1.Range("F7").Value = 33
$ws1.Range("G7").Value = 33
$ws1.Range("H7").Value = 19
$ws1.Range("I7").Value = 11

$ws1.Range("A8").Value = "Students More than 75%"
$ws1.Range("B8").Value = 0
$ws1.Range("C8").Value = 5
$ws1.Range("D8").Value = 13
$ws1.Range("E8").Value = 23
$ws1.Range("F8").Value = 7
$ws1.Range("G8").Value = 3
$ws1.Range("H8").Value = 25
$ws1.Range("I8").Value = 28

# ---- Sheet 2: Slow Learners ----
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A1").Value = "Roll No."
$ws2.Range("B1").Value = "Student Name"
$ws2.Range("C1").Value = "Count"

$ws2.Range("A2").Value = 2236150002
$ws2.Range("B2").Value = "Mohammed Azees M"
$ws2.Range("C2").Value = 7

$ws2.Range("A3").Value = 2136110035
$ws2.Range("B3").Value = "Manuneethi S"
$ws2.Range("C3").Value = 7

$ws2.Range("A4").Value = 2136110029
$ws2.Range("B4").Value = "Balaganapathi A"
$ws2.Range("C4").Value = 7

$ws2.Range("A5").Value = 2236150003
$ws2.Range("B5").Value = "Krishnakumar S"
$ws2.Range("C5").Value = 6

$ws2.Range("A6").Value = 2136110046
$ws2.Range("B6").Value = "MARIKANNAN P"
$ws2.Range("C6").Value = 5

$ws2.Range("A7").Value = 2136110044
$ws2.Range("B7").Value = "SATHISHKUMAR N"
$ws2.Range("C7").Value = 5

$ws2.Range("A8").Value = 2136110024
$ws2.Range("B8").Value = "Upanshu"
$ws2.Range("C8").Value = 5

$ws2.Range("A9").Value = 2136110033
$ws2.Range("B9").Value = "JAISANKAR S"
$ws2.Range("C9").Value = 5

$ws2.Range("A10").Value = 2136110043
$ws2.Range("B10").Value = "Sanjay S"
$ws2.Range("C10").Value = 4

$ws2.Range("A11").Value = 2136110006
$ws2.Range("B11").Value = "DHINAKARAN R"
$ws2.Range("C11").Value = 4

$ws2.Range("A12").Value = 2136110011
$ws2.Range("B12").Value = "Manoharan K"
$ws2.Range("C12").Value = 4

$ws2.Range("A13").Value = 2136110002
$ws2.Range("B13").Value = "Arulselvam C"
$ws2.Range("C13").Value = 4

$ws2.Range("A14").Value = 2136110032
$ws2.Range("B14").Value = "Jaikrishnan V"
$ws2.Range("C14").Value = 4

# ---- Sheet 3: Fast Learners ----
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A1").Value = "Roll No."
$ws3.Range("B1").Value = "Student Name"
$ws3.Range("C1").Value = "Count"

$ws3.Range("A2").Value = 2136110001
$ws3.Range("B2").Value = "Aravind S"
$ws3.Range("C2").Value = 7

$ws3.Range("A3").Value = 2136110031
$ws3.Range("B3").Value = "Hitesh Kumar K A"
$ws3.Range("C3").Value = 6

$ws3.Range("A4").Value = 2136110013
$ws3.Range("B4").Value = "Naveena A"
$ws3.Range("C4").Value = 6

$ws3.Range("A5").Value = 2136110008
$ws3.Range("B5").Value = "Jananika B"
$ws3.Range("C5").Value = 6

$ws3.Range("A6").Value = 2136110019
$ws3.Range("B6").Value = "Sivaa Ganesh S"
$ws3.Range("C6").Value = 6

$ws3.Range("A7").Value = 2136110030
$ws3.Range("B7").Value = "Brijesh A"
$ws3.Range("C7").Value = 6

$ws3.Range("A8").Value = 2136110003
$ws3.Range("B8").Value = "Ashik Jenly V L"
$ws3.Range("C8").Value = 5

$ws3.Range("A9").Value = 2136110021
$ws3.Range("B9").Value = "Subhashini S"
$ws3.Range("C9").Value = 5

$ws3.Range("A10").Value = 2136110016
$ws3.Range("B10").Value = "Nithya Sri R"
$ws3.Range("C10").Value = 5

$ws3.Range("A11").Value = 2136110040
$ws3.Range("B11").Value = "Preethiga S"
$ws3.Range("C11").Value = 5

$ws3.Range("A12").Value = 2136110009
$ws3.Range("B12").Value = "Kalaivani S"
$ws3.Range("C12").Value = 5

$ws3.Range("A13").Value = 2136110026
$ws3.Range("B13").Value = "AJAY S"
$ws3.Range("C13").Value = 5

$ws3.Range("A14").Value = 2136110022
$ws3.Range("B14").Value = "Suji Shri B"
$ws3.Range("C14").Value = 5

$ws3.Range("A15").Value = 2136110038
$ws3.Range("B15").Value = "Pradeep M"
$ws3.Range("C15").Value = 5

$ws3.Range("A16").Value = 2136110047
$ws3.Range("B16").Value = "Gowtham R"
$ws3.Range("C16").Value = 5

$ws3.Range("A17").Value = 2136110049
$ws3.Range("B17").Value = "Kailashwaran R"
$ws3.Range("C17").Value = 5

$ws3.Range("A18").Value = 2136110036
$ws3.Range("B18").Value = "Mohamed Suhail J"
$ws3.Range("C18").Value = 4

$ws3.Range("A19").Value = 2136110010
$ws3.Range("B19").Value = "Krishnapriya K"
$ws3.Range("C19").Value = 4

$ws3.Range("A20").Value = 2136110004
$ws3.Range("B20").Value = "Deepakragavan J"
$ws3.Range("C20").Value = 4

$ws3.Range("A21").Value = 2136110007
$ws3.Range("B21").Value = "Guruprasath V"
$ws3.Range("C21").Value = 4

$ws3.Range("A22").Value = 2136110045
$ws3.Range("B22").Value = "Varsha V"
$ws3.Range("C22").Value = 4

$ws3.Range("A23").Value = 2136110014
$ws3.Range("B23").Value = "Nawin B"
$ws3.Range("C23").Value = 4

$ws3.Range("A24").Value = 2136110020
$ws3.Range("B24").Value = "Srija D"
$ws3.Range("C24").Value = 4
